$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Laura / New Denizen
$ws.Range("A8").Value = "Laura"
$ws.Range("B8").Value = "New Denizen"
$ws.Range("C8").Value = "newdenizenblog@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:newdenizenblog@gmail.com")
$ws.Range("D8").Value = "food blogger"

# Row 9 - Josie Sexton / Eater Denver
$ws.Range("A9").Value = "Josie Sexton"
$ws.Range("B9").Value = "Eater Denver"
$ws.Range("C9").Value = "https://www.josiesexton.com/contact.html"
$ws.Range("D9").Value = "Food blogger at multiple organizations (was at the Coloradoan in the past?)"

# Row 10 - Mile High and Hungry
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:milehighandhungry@gmail.com", "", "", "mailto:milehighandhungry@gmail.com")
$ws.Range("B10").Value = "milehighandhungry@gmail.com"
$ws.Range("C10").Value = "http://milehighandhungry.com/"
$ws.Range("D10").Value = "Instagram influencer and blogger"
$ws.Range("E10").Value = "They have a favorite hh spot"

# Row 11 - Best Booze Denver (instagram)
$ws.Range("C11").Value = "https://www.instagram.com/bestboozedenver/"
$ws.Range("D11").Value = "Instagram influencer and blogger"

# Row 12 - Denver food scene (instagram)
$ws.Range("C12").Value = "https://www.instagram.com/Denverfoodscene/"
$ws.Range("H12").Value = "Has already posted about an app"

# Row 13 - Cara Chancellor / 303 magazine
$ws.Range("A13").Value = "Cara Chancellor"
$ws.Range("B13").Value = "cara.chancellor@yahoo.com"
$ws.Range("C13").Value = "303 magazine"

# Column widths to roughly match target best-fit widths
$ws.Columns("A:A").ColumnWidth = 14
$ws.Columns("B:B").ColumnWidth = 29.17
$ws.Columns("C:C").ColumnWidth = 43.83
$ws.Columns("D:D").ColumnWidth = 68.33

# Restore selection similar to the original authoring state
$ws.Range("C14").Select()
